$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$values = @{
    "H12" = 361
    "J12" = 0
    "L12" = 0
    "H33" = 416926.75
    "I33" = 625227.75
    "J33" = 324.75
    "K33" = 625227.75
    "L33" = 324.75
    "M33" = -624998.75
    "N33" = -782.75
    "H41" = 1479.7273
    "I41" = 896.8570999999999
    "K41" = 896.8570999999999
    "M41" = -456.8570999999999
    "H86" = 3590203.5
    "I86" = 6458252
    "J86" = 5143.25
    "K86" = 6458252
    "L86" = 5143.25
    "M86" = -6457129
    "N86" = -7389.25
    "H89" = 3590203.5
    "I89" = 6458252
    "J89" = 5143.25
    "K89" = 32291260
    "L89" = 25716.25
    "M89" = -32285644
    "N89" = -36948.25
    "H101" = 3021.889
    "I101" = 601.6
    "J101" = 6047.25
    "K101" = 1804.8
    "L101" = 18141.75
    "M101" = -182.8000000000002
    "N101" = -21385.75
    "H103" = 1495
    "I103" = 1495
    "K103" = 4485
    "M103" = -3899
    "H111" = 2606
    "I111" = 2606
    "K111" = 7818
    "M111" = -4751
    "H116" = 41054984
    "I116" = 38026844
    "J116" = 45459550
    "K116" = 38026844
    "L116" = 45459550
    "M116" = -38023402
    "N116" = -45466434
    "H121" = 3762.1155
    "J121" = 3762.1155
    "L121" = 11286.3465
    "N121" = -14780.3465
    "H127" = 2694.9443
    "I127" = 1399.6923
    "J127" = 3427.0435
    "K127" = 4199.0769
    "L127" = 10281.1305
    "M127" = 760.9231
    "N127" = -20201.1305
    "H132" = 11015.655
    "I132" = 1814.5555
    "J132" = 14867.279
    "K132" = 5443.666499999999
    "L132" = 44601.837
    "M132" = -2913.666499999999
    "N132" = -49661.837
    "H135" = 2091.2964
    "I135" = 427.8421
    "J135" = 6042
    "K135" = 3850.5789
    "L135" = 54378
    "M135" = -1315.5789
    "N135" = -59448
    "H137" = 6805723
    "I137" = 1159.8572
    "K137" = 3479.5716
    "M137" = -929.5715999999998
    "H141" = 5956.08
    "I141" = 4723.091
    "K141" = 14169.273
    "M141" = -8989.273000000001
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
$ws.Range("N12").ClearContents()

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$values = @{
    "H32" = 5002.1177
    "I32" = 3018.1082
    "K32" = 3018.1082
    "M32" = -2731.1082
    "H45" = 5000
    "I45" = 6000
    "K45" = 6000
    "M45" = -5623
    "H46" = 10000
    "J46" = 10000
    "L46" = 10000
    "N46" = -10638
    "H51" = 77071.336
    "J51" = 77071.336
    "L51" = 77071.336
    "N51" = -78583.336
    "H61" = 20108.54
    "I61" = 24514.572
    "J61" = 14968.167
    "K61" = 24514.572
    "L61" = 14968.167
    "M61" = -24302.572
    "N61" = -15392.167
    "H102" = 572223.25
    "I102" = 1371191.6
    "K102" = 1371191.6
    "M102" = -1369569.6
    "H132" = 17192.545
    "I132" = 20496.8
    "K132" = 61490.39999999999
    "M132" = -58960.39999999999
    "H136" = 20108.54
    "I136" = 24514.572
    "J136" = 14968.167
    "K136" = 73543.716
    "L136" = 44904.501
    "M136" = -70993.716
    "N136" = -50004.501
    "H139" = 95374.5
    "I139" = 174968
    "K139" = 174968
    "M139" = -169828
    "H140" = 114650
    "J140" = 114650
    "L140" = 114650
    "N140" = -125010
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$values = @{
    "H80" = 334.45456
    "J80" = 217.375
    "L80" = 217.375
    "N80" = -2213.375
    "H83" = 334.45456
    "J83" = 217.375
    "L83" = 1086.875
    "N83" = -11070.875
    "H94" = 914251.75
    "I94" = 4566925.5
    "J94" = 1083.25
    "K94" = 4566925.5
    "L94" = 1083.25
    "M94" = -4566474.5
    "N94" = -1985.25
    "H99" = 2316625.2
    "I99" = 4167965.5
    "J99" = 2450
    "K99" = 4167965.5
    "L99" = 2450
    "M99" = -4166467.5
    "N99" = -5446
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$values = @{
    "H17" = 10000
    "I17" = 10000
    "K17" = 10000
    "M17" = -9826
    "H31" = 3032.0715
    "J31" = 5367.5835
    "L31" = 5367.5835
    "N31" = -5957.5835
    "H34" = 3032.0715
    "J34" = 5367.5835
    "L34" = 5367.5835
    "N34" = -5771.5835
    "H105" = 1059.6
    "I105" = 893.7646999999999
    "J105" = 1999.3334
    "K105" = 893.7646999999999
    "L105" = 1999.3334
    "M105" = 853.2353000000001
    "N105" = -5493.3334
    "H134" = 2408.0527
    "I134" = 2338
    "K134" = 7014
    "M134" = -4479
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$values = @{
    "H96" = 12996.5
    "I96" = 12986
    "J96" = 13000
    "K96" = 38958
    "L96" = 39000
    "M96" = -36899
    "N96" = -43118
    "H137" = 6123047.5
    "I137" = 5275
    "J137" = 7346602
    "K137" = 15825
    "L137" = 22039806
    "M137" = -10725
    "N137" = -22050006
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$values = @{
    "H113" = 3803.6667
    "I113" = 3705.5
    "K113" = 3705.5
    "M113" = -1535.5
    "H126" = 4769.75
    "I126" = 2780.7693
    "J126" = 8463.571
    "K126" = 8342.3079
    "L126" = 25390.713
    "M126" = -5872.3079
    "N126" = -30330.713
    "H138" = 60000
    "I138" = 0
    "K138" = 0
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
$ws.Range("M138").ClearContents()

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$values = @{
    "H5" = 5000
    "J5" = 5000
    "L5" = 5000
    "N5" = -5226
    "H31" = 14937.5
    "J31" = 14937.5
    "L31" = 14937.5
    "N31" = -15433.5
    "H40" = 6989.75
    "I40" = 6601.4
    "K40" = 6601.4
    "M40" = -6465.4
    "H46" = 6283.467
    "I46" = 3858.1667
    "J46" = 6889.7915
    "K46" = 3858.1667
    "L46" = 6889.7915
    "M46" = -3670.1667
    "N46" = -7265.7915
    "H100" = 3490.889
    "I100" = 1188.5
    "K100" = 1188.5
    "M100" = -647.5
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$values = @{
    "H43" = 79998
    "J43" = 79998
    "L43" = 79998
    "N43" = -80296
    "H55" = 21009.334
    "J55" = 21009.334
    "L55" = 21009.334
    "N55" = -21563.334
    "H107" = 2624.923
    "I107" = 2928.9048
    "J107" = 1348.2
    "K107" = 8786.714399999999
    "L107" = 4044.6
    "M107" = -6866.714399999999
    "N107" = -7884.6
    "H136" = 8759.5
    "J136" = 9309.690000000001
    "L136" = 27929.07
    "N136" = -33029.07
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
